# Update "想去人数" (want-to-go count) / "最低票价" (lowest price) figures
# across the four sheets, matching the refreshed scrape output at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 367
$ws.Range("G6").Value = 80
$ws.Range("F8").Value = 809
$ws.Range("F9").Value = 4205
$ws.Range("F11").Value = 175
$ws.Range("F13").Value = 6082
$ws.Range("F16").Value = 2335
$ws.Range("F19").Value = 472
$ws.Range("F20").Value = 9152
$ws.Range("F22").Value = 2458
$ws.Range("F24").Value = 2311
$ws.Range("F25").Value = 2437
$ws.Range("F26").Value = 1392
$ws.Range("F28").Value = 1960
$ws.Range("F30").Value = 58
$ws.Range("F31").Value = 330
$ws.Range("F34").Value = 281
$ws.Range("F39").Value = 1216
$ws.Range("F42").Value = 238
$ws.Range("F43").Value = 1529
$ws.Range("F44").Value = 2513
$ws.Range("F45").Value = 924
$ws.Range("F46").Value = 295
$ws.Range("F48").Value = 20
$ws.Range("F49").Value = 24

# --- Sheet "演出" (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F22").Value = 67
$ws.Range("F23").Value = 67

# --- Sheet "本地生活" (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 891

# --- Sheet "全部类型" (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 891
$ws.Range("F7").Value = 367
$ws.Range("G9").Value = 80
$ws.Range("F13").Value = 809
$ws.Range("F14").Value = 4205
$ws.Range("F15").Value = 175
$ws.Range("F16").Value = 6082
$ws.Range("F19").Value = 2335
$ws.Range("F22").Value = 9152
$ws.Range("F24").Value = 2458
$ws.Range("F25").Value = 2311
$ws.Range("F26").Value = 1392
$ws.Range("F28").Value = 1960
$ws.Range("F30").Value = 58
$ws.Range("F31").Value = 330
$ws.Range("F33").Value = 281
$ws.Range("F37").Value = 1216
$ws.Range("F40").Value = 238
$ws.Range("F41").Value = 1529
$ws.Range("F42").Value = 2513
$ws.Range("F43").Value = 924
$ws.Range("F44").Value = 295
$ws.Range("F49").Value = 20
$ws.Range("F50").Value = 67
